$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Ptn", "Ptprs", "ECs", [double]"3", [double]"1", [double]"2.327816333333333", [double]"6.983449", [double]"0.02128501190197005", [double]"0.02128501190197004", [double]"3", [double]"1", [double]"4.925988333333333", [double]"14.777965", [double]"0.05656988822582037", [double]"0.05656988822582035", [double]"11.46679610014278", [double]"103.201164901285", [double]"0.001204090744179702", [double]"0.001204090744179701"),
  @("ECs", "Ptn", "Ptprs", "FAPs", [double]"3", [double]"1", [double]"2.327816333333333", [double]"6.983449", [double]"0.02128501190197005", [double]"0.02128501190197004", [double]"3", [double]"1", [double]"47.585289", [double]"142.755867", [double]"0.5464678959362861", [double]"0.5464678959362861", [double]"110.769812960587", [double]"996.928316645283", [double]"0.01163157566904838", [double]"0.01163157566904838"),
  @("ECs", "Ptn", "Ptprs", "Inflammatory-Mac", [double]"3", [double]"1", [double]"2.327816333333333", [double]"6.983449", [double]"0.02128501190197005", [double]"0.02128501190197004", [double]"3", [double]"1", [double]"13.80191933333334", [double]"41.40575800000001", [double]"0.1585007882996995", [double]"0.1585007882996994", [double]"32.12833325548245", [double]"289.154999299342", [double]"0.003373691165430738", [double]"0.003373691165430737"),
  @("ECs", "Ptn", "Ptprs", "MuSCs", [double]"3", [double]"1", [double]"2.327816333333333", [double]"6.983449", [double]"0.02128501190197005", [double]"0.02128501190197004", [double]"3", [double]"1", [double]"12.18796133333333", [double]"36.563884", [double]"0.1399661476381804", [double]"0.1399661476381803", [double]"28.37133546176845", [double]"255.342019155916", [double]"0.002979181118351566", [double]"0.002979181118351565"),
  @("ECs", "Ptn", "Ptprs", "Resolving-Mac", [double]"3", [double]"1", [double]"2.327816333333333", [double]"6.983449", [double]"0.02128501190197005", [double]"0.02128501190197004", [double]"3", [double]"1", [double]"8.576764333333333", [double]"25.730293", [double]"0.09849527990001386", [double]"0.09849527990001385", [double]"19.96513210228411", [double]"179.686188920557", [double]"0.002096473204959666", [double]"0.002096473204959665"),
  @("FAPs", "Ptn", "Ptprs", "ECs", [double]"3", [double]"1", [double]"103.0385286666667", [double]"309.115586", [double]"0.9421603747796319", [double]"0.9421603747796318", [double]"3", [double]"1", [double]"4.925988333333333", [double]"14.777965", [double]"0.05656988822582037", [double]"0.05656988822582035", [double]"507.5665900958322", [double]"4568.09931086249", [double]"0.0532979070920808", [double]"0.05329790709208079"),
  @("FAPs", "Ptn", "Ptprs", "FAPs", [double]"3", [double]"1", [double]"103.0385286666667", [double]"309.115586", [double]"0.9421603747796319", [double]"0.9421603747796318", [double]"3", [double]"1", [double]"47.585289", [double]"142.755867", [double]"0.5464678959362861", [double]"0.5464678959362861", [double]"4903.118164738117", [double]"44128.06348264306", [double]"0.5148603976403682", [double]"0.5148603976403681"),
  @("FAPs", "Ptn", "Ptprs", "Inflammatory-Mac", [double]"3", [double]"1", [double]"103.0385286666667", [double]"309.115586", [double]"0.9421603747796319", [double]"0.9421603747796318", [double]"3", [double]"1", [double]"13.80191933333334", [double]"41.40575800000001", [double]"0.1585007882996995", [double]"0.1585007882996994", [double]"1422.129460882688", [double]"12799.16514794419", [double]"0.1493331621073119", [double]"0.1493331621073119"),
  @("FAPs", "Ptn", "Ptprs", "MuSCs", [double]"3", [double]"1", [double]"103.0385286666667", [double]"309.115586", [double]"0.9421603747796319", [double]"0.9421603747796318", [double]"3", [double]"1", [double]"12.18796133333333", [double]"36.563884", [double]"0.1399661476381804", [double]"0.1399661476381803", [double]"1255.829603232891", [double]"11302.46642909602", [double]"0.1318705581152493", [double]"0.1318705581152493"),
  @("FAPs", "Ptn", "Ptprs", "Resolving-Mac", [double]"3", [double]"1", [double]"103.0385286666667", [double]"309.115586", [double]"0.9421603747796319", [double]"0.9421603747796318", [double]"3", [double]"1", [double]"8.576764333333333", [double]"25.730293", [double]"0.09849527990001386", [double]"0.09849527990001385", [double]"883.7371776274108", [double]"7953.634598646698", [double]"0.0927983498246218", [double]"0.09279834982462178"),
  @("MuSCs", "Ptn", "Ptprs", "ECs", [double]"3", [double]"1", [double]"3.975769333333333", [double]"11.927308", [double]"0.03635351138648862", [double]"0.03635351138648861", [double]"3", [double]"1", [double]"4.925988333333333", [double]"14.777965", [double]"0.05656988822582037", [double]"0.05656988822582035", [double]"19.58459335202444", [double]"176.26134016822", [double]"0.002056514075749749", [double]"0.002056514075749748"),
  @("MuSCs", "Ptn", "Ptprs", "FAPs", [double]"3", [double]"1", [double]"3.975769333333333", [double]"11.927308", [double]"0.03635351138648862", [double]"0.03635351138648861", [double]"3", [double]"1", [double]"47.585289", [double]"142.755867", [double]"0.5464678959362861", [double]"0.5464678959362861", [double]"189.188132724004", [double]"1702.693194516036", [double]"0.01986602687727025", [double]"0.01986602687727025"),
  @("MuSCs", "Ptn", "Ptprs", "Inflammatory-Mac", [double]"3", [double]"1", [double]"3.975769333333333", [double]"11.927308", [double]"0.03635351138648862", [double]"0.03635351138648861", [double]"3", [double]"1", [double]"13.80191933333334", [double]"41.40575800000001", [double]"0.1585007882996995", [double]"0.1585007882996994", [double]"54.87324762660712", [double]"493.859228639464", [double]"0.005762060212220546", [double]"0.005762060212220544"),
  @("MuSCs", "Ptn", "Ptprs", "MuSCs", [double]"3", [double]"1", [double]"3.975769333333333", [double]"11.927308", [double]"0.03635351138648862", [double]"0.03635351138648861", [double]"3", [double]"1", [double]"12.18796133333333", [double]"36.563884", [double]"0.1399661476381804", [double]"0.1399661476381803", [double]"48.45652290491911", [double]"436.108706144272", [double]"0.005088260941887537", [double]"0.005088260941887535"),
  @("MuSCs", "Ptn", "Ptprs", "Resolving-Mac", [double]"3", [double]"1", [double]"3.975769333333333", [double]"11.927308", [double]"0.03635351138648862", [double]"0.03635351138648861", [double]"3", [double]"1", [double]"8.576764333333333", [double]"25.730293", [double]"0.09849527990001386", [double]"0.09849527990001385", [double]"34.09923661569378", [double]"306.893129541244", [double]"0.003580649279360537", [double]"0.003580649279360536"),
  @("Resolving-Mac", "Ptn", "Ptprs", "ECs", [double]"1", [double]"0.3333333333333333", [double]"0.02199333333333333", [double]"0.06598", [double]"0.0002011019319095741", [double]"0.0002011019319095741", [double]"3", [double]"1", [double]"4.925988333333333", [double]"14.777965", [double]"0.05656988822582037", [double]"0.05656988822582035", [double]"0.1083389034111111", [double]"0.9750501306999998", [double]"1.137631381012115E-05", [double]"1.137631381012114E-05"),
  @("Resolving-Mac", "Ptn", "Ptprs", "FAPs", [double]"1", [double]"0.3333333333333333", [double]"0.02199333333333333", [double]"0.06598", [double]"0.0002011019319095741", [double]"0.0002011019319095741", [double]"3", [double]"1", [double]"47.585289", [double]"142.755867", [double]"0.5464678959362861", [double]"0.5464678959362861", [double]"1.04655912274", [double]"9.419032104659999", [double]"0.0001098957495993473", [double]"0.0001098957495993472"),
  @("Resolving-Mac", "Ptn", "Ptprs", "Inflammatory-Mac", [double]"1", [double]"0.3333333333333333", [double]"0.02199333333333333", [double]"0.06598", [double]"0.0002011019319095741", [double]"0.0002011019319095741", [double]"3", [double]"1", [double]"13.80191933333334", [double]"41.40575800000001", [double]"0.1585007882996995", [double]"0.1585007882996994", [double]"0.3035502125377778", [double]"2.73195191284", [double]"3.187481473625998E-05", [double]"3.187481473625998E-05"),
  @("Resolving-Mac", "Ptn", "Ptprs", "MuSCs", [double]"1", [double]"0.3333333333333333", [double]"0.02199333333333333", [double]"0.06598", [double]"0.0002011019319095741", [double]"0.0002011019319095741", [double]"3", [double]"1", [double]"12.18796133333333", [double]"36.563884", [double]"0.1399661476381804", [double]"0.1399661476381803", [double]"0.2680538962577778", [double]"2.41248506632", [double]"2.814746269197875E-05", [double]"2.814746269197874E-05"),
  @("Resolving-Mac", "Ptn", "Ptprs", "Resolving-Mac", [double]"1", [double]"0.3333333333333333", [double]"0.02199333333333333", [double]"0.06598", [double]"0.0002011019319095741", [double]"0.0002011019319095741", [double]"3", [double]"1", [double]"8.576764333333333", [double]"25.730293", [double]"0.09849527990001386", [double]"0.09849527990001385", [double]"0.1886316369044445", [double]"1.69768473214", [double]"1.980759107186703E-05", [double]"1.980759107186703E-05")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}

Write-Output $ws.Range("A2").Value()
Write-Output $ws.Range("D21").Value()
Write-Output $ws.Range("T21").Value()
